$wb = $excel.ActiveWorkbook

# New metric values (B2:B14) for each affected worksheet.
$updates = @{
    "Bidirectional A" = @(
        [double]"0.0001339589998679003",
        [double]"0.0001431670007150387",
        [double]"1600.8",
        [double]"1600.8",
        [double]"0.00240325927734375",
        [double]"0.00233612060546875",
        [double]"0.005446434020996094",
        [double]"0.005446434020996094",
        [double]"2.541700087022036e-05",
        [double]"3.162500070175156e-05",
        [double]"3.118739987257868e-05",
        [double]"2.790010021271883e-05",
        [double]"10"
    )
    "D Lite" = @(
        [double]"0.0001175840006908402",
        [double]"0.0006950419992790557",
        [double]"1600.8",
        [double]"1600.8",
        [double]"0.0023193359375",
        [double]"0.0023193359375",
        [double]"0.07711029052734375",
        [double]"0.07694091796874999",
        [double]"0.0001895830009743804",
        [double]"0.0002845850031008013",
        [double]"2.869989966711728e-05",
        [double]"2.476270001352532e-05",
        [double]"10"
    )
    "IDA" = @(
        [double]"0.0001976250005100155",
        [double]"0.000957290998485405",
        [double]"1600.8",
        [double]"1600.8",
        [double]"0.0023193359375",
        [double]"0.0023193359375",
        [double]"0.002044677734375",
        [double]"0.00198974609375",
        [double]"3.095799911534414e-05",
        [double]"0.0003507499986881157",
        [double]"0.00241379999970377",
        [double]"3.574999955162639e-05",
        [double]"10"
    )
    "SMA" = @(
        [double]"0.0001421249999111751",
        [double]"7.666599958611187e-05",
        [double]"1600.8",
        [double]"1600.8",
        [double]"0.0023193359375",
        [double]"0.0023193359375",
        [double]"0.00255584716796875",
        [double]"0.00255584716796875",
        [double]"3.099999958067201e-05",
        [double]"6.374999975378159e-05",
        [double]"7.805820023349952e-05",
        [double]"3.347500023664907e-05",
        [double]"10"
    )
    "RTAA (L=25, M=3)" = @(
        [double]"0.0001369590008835075",
        [double]"0.0001626670000405284",
        [double]"1600.8",
        [double]"1600.8",
        [double]"0.0023193359375",
        [double]"0.0023193359375",
        [double]"0.00458526611328125",
        [double]"0.00458526611328125",
        [double]"2.791699989757035e-05",
        [double]"6.650000068475492e-05",
        [double]"6.821679980930639e-05",
        [double]"2.978729971800931e-05",
        [double]"10"
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $values = $updates[$sheetName]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $row = $i + 2
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}
